$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data for rows 2-13 (columns A-T), reflecting the new TPM-based
# recalculation plus the addition of the "Inflammatory-Mac" target cluster.
$data = @(
    @("FAPs","Inha","Tgfbr3","ECs",2,0.6666666666666666,0.3688166666666666,1.10645,0.6243561829432177,0.7137249741492807,2,1,22.3363875,44.672775,0.2333075346730891,0.1773553735433322,8.238031983125,49.42819189875,0.1456670018003823,0.1265829593974508),
    @("FAPs","Inha","Tgfbr3","FAPs",2,0.6666666666666666,0.3688166666666666,1.10645,0.6243561829432177,0.7137249741492807,3,1,60.113367,180.340101,0.627894795237978,0.715968192656428,22.17081163905,199.53730475145,0.3920299976446972,0.5110043797954162),
    @("FAPs","Inha","Tgfbr3","Inflammatory-Mac",2,0.6666666666666666,0.3688166666666666,1.10645,0.6243561829432177,0.7137249741492807,3,1,0.1245063333333333,0.373519,0.001300490765625635,0.001482907694239543,0.04592001083888888,0.4132800975499999,0.0008119694503789241,0.001058388255736887),
    @("FAPs","Inha","Tgfbr3","MuSCs",2,0.6666666666666666,0.3688166666666666,1.10645,0.6243561829432177,0.7137249741492807,2,1,12.9946545,25.989309,0.1357314742736964,0.103180149561519,4.792645157175,28.75587094305,0.08474478518278061,0.07364224957851405),
    @("FAPs","Inha","Tgfbr3","Neutrophils",2,0.6666666666666666,0.3688166666666666,1.10645,0.6243561829432177,0.7137249741492807,2,0.6666666666666666,0.121553,0.364659,0.001269642674408205,0.00144773261031888,0.04483077228333333,0.40347695055,0.0007927092538953255,0.001033282919874913),
    @("FAPs","Inha","Tgfbr3","Resolving-Mac",2,0.6666666666666666,0.3688166666666666,1.10645,0.6243561829432177,0.7137249741492807,3,1,0.047492,0.142476,0.000496062375202541,0.0005656439341625811,0.01751584113333333,0.1576425702,0.0003097196110832048,0.0004037142022878856),
    @("FAPs","Inha","Tgfbr3","ECs",1,0.5,0.2218985,0.443797,0.3756438170567823,0.2862750258507193,2,1,22.3363875,44.672775,0.2333075346730891,0.1773553735433322,4.95641088166875,19.825643526675,0.08764053287270679,0.05077241414588141),
    @("FAPs","Inha","Tgfbr3","FAPs",1,0.5,0.2218985,0.443797,0.3756438170567823,0.2862750258507193,3,1,60.113367,180.340101,0.627894795237978,0.715968192656428,13.3390659672495,80.034395803497,0.2358647975932808,0.2049638128610117),
    @("FAPs","Inha","Tgfbr3","Inflammatory-Mac",1,0.5,0.2218985,0.443797,0.3756438170567823,0.2862750258507193,3,1,0.1245063333333333,0.373519,0.001300490765625635,0.001482907694239543,0.02762776860716666,0.165766611643,0.0004885213152467108,0.0004245194385026557),
    @("FAPs","Inha","Tgfbr3","MuSCs",1,0.5,0.2218985,0.443797,0.3756438170567823,0.2862750258507193,2,1,12.9946545,25.989309,0.1357314742736964,0.103180149561519,2.88349434156825,11.533977366273,0.05098668909091576,0.02953789998300493),
    @("MuSCs","Inha","Tgfbr3","Neutrophils",1,0.5,0.2218985,0.443797,0.3756438170567823,0.2862750258507193,2,0.6666666666666666,0.121553,0.364659,0.001269642674408205,0.00144773261031888,0.0269724283705,0.161834570223,0.0004769334205128797,0.0004144496904439666),
    @("MuSCs","Inha","Tgfbr3","Resolving-Mac",1,0.5,0.2218985,0.443797,0.3756438170567823,0.2862750258507193,3,1,0.047492,0.142476,0.000496062375202541,0.0005656439341625811,0.010538403562,0.063230421372,0.0001863427641193363,0.0001619297318746955)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Count; $j++) {
        $ws.Cells.Item($r, $j + 1).Value = $row[$j]
    }
}

$wb.Save()
